$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.103302717208862
$ws.Range("B1").Value = 2.633480548858643
$ws.Range("C1").Value = 2.02814245223999
$ws.Range("D1").Value = 1.819197773933411
$ws.Range("E1").Value = 1.830005526542664
